$wb = $excel.ActiveWorkbook

# --- Rename "Sheet3" to "SampleRegistration" (sheetId 26) ---
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "SampleRegistration"

# --- Populate the new SampleRegistration sheet with the Load Test Suite data ---
$ws.Range("A1").Value = "Product_Code"
$ws.Range("B1").Value = "Specification_Name"
$ws.Range("A2").Value = "Auto_Prod_1"
$ws.Range("B2").Value = "Auto_Spec_1"
$ws.Range("A3").Value = "Auto_Prod_2"
$ws.Range("B3").Value = "Auto_Spec_2"

# Apply thin black "all borders" box to the populated range
$ws.Range("A1:B3").Borders.ColorIndex = 1
$ws.Range("A1:B3").Borders.LineStyle = 1

# Column widths (calibrated so stored OOXML width matches 14 / 19)
$ws.Columns.Item(1).ColumnWidth = 13.166666666666666
$ws.Columns.Item(2).ColumnWidth = 18.166666666666668

# Page setup (portrait, A4 paper)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

$ws.Range("E14").Select()

# --- Touch UserCreation rows 3:9 (re-entered without border formatting) ---
$uc = $wb.Worksheets.Item("UserCreation")
$uc.Range("A3:F9").ClearFormats()
$uc.Activate()
$uc.Range("D14").Select()
